$d = $word.ActiveDocument

# The document currently has two paragraphs:
#   1) "Async in Rust:" (bold title)
#   2) an empty paragraph (holds the _GoBack bookmark) in normal (size 24) text
#
# Insert a brand new paragraph *before* the trailing empty paragraph so the
# new paragraph inherits that paragraph's (non-bold, size-24) formatting
# rather than the bold title formatting.

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs(2)
$newPara.Range.Text = "Async program is a form of parallel programming, its allows the atomic work to run separately from the primary (main) application thread. When the atomic work is complete, it notifies the main thread with failed or success status. Moreover, Asynchronous code allows us to run multiple tasks concurrently on the same OS thread."
